$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "is_public" column (O) first, so its header string is registered
#     in the shared-strings table before "project_id" -- matches the
#     original authoring order captured in the target file.
$ws.Range("O1").Value2 = "is_public"
$ws.Range("O2").Value2 = $true
$ws.Range("O3").Value2 = $true
$ws.Range("O4").Value2 = $true
$ws.Range("O5").Value2 = $true
$ws.Range("O6").Value2 = $true
$ws.Range("O7").Value2 = $true
$ws.Range("O8").Value2 = $true
$ws.Range("O9").Value2 = $true

# --- Add "project_id" column (N), populated per the project each row
#     belongs to (matches B/C/L columns already on each row).
$ws.Range("N1").Value2 = "project_id"
$ws.Range("N2").Value2 = "neear"
$ws.Range("N3").Value2 = "neear"
$ws.Range("N4").Value2 = "neear"
$ws.Range("N5").Value2 = "aia_crafting_the_future"
$ws.Range("N6").Value2 = "aia_crafting_the_future"
$ws.Range("N7").Value2 = "aia_crafting_the_future"
$ws.Range("N8").Value2 = "ajo_artisan_lofts"
$ws.Range("N9").Value2 = "building_with_bamboo_desert"

# Column M (project "About" tags) ends up auto-fit to its widest entry
$ws.Columns.Item(13).ColumnWidth = 18.998697916666668

# Final cursor position left by the editor after keying in the last value
$ws.Range("O16").Select() | Out-Null
